$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 5.1
$ws.Range("N2").Value = 3.4
$ws.Range("P2").Value = 1.8
$ws.Range("T2").Value = 1.99
$ws.Range("V2").Value = 1.24
$ws.Range("Y2").Value = 15.5
$ws.Range("AI2").Value = 1000
$ws.Range("AL2").Value = 44
$ws.Range("N3").Value = 3.35
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.8
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 4.1
$ws.Range("U3").Value = 1.9
$ws.Range("Y3").Value = 7.6
$ws.Range("AA3").Value = 18
$ws.Range("AC3").Value = 8.199999999999999
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 44
$ws.Range("AJ3").Value = 150
$ws.Range("AK3").Value = 85
$ws.Range("AL3").Value = 90
$ws.Range("AM3").Value = 150
$ws.Range("H4").Value = 1.84
$ws.Range("P4").Value = 1.65
$ws.Range("R4").Value = 1.24
$ws.Range("U4").Value = 1.75
$ws.Range("W4").Value = 1.21
$ws.Range("Y4").Value = 6.6
$ws.Range("F5").Value = 5.1
$ws.Range("G5").Value = 5.3
$ws.Range("H5").Value = 1.8
$ws.Range("I5").Value = 1.82
$ws.Range("J5").Value = 3.9
$ws.Range("L5").Value = 1.4
$ws.Range("P5").Value = 1.9
$ws.Range("Q5").Value = 2.06
$ws.Range("R5").Value = 1.34
$ws.Range("S5").Value = 3.75
$ws.Range("V5").Value = 2.22
$ws.Range("W5").Value = 1.23
$ws.Range("X5").Value = 13
$ws.Range("Y5").Value = 8
$ws.Range("AA5").Value = 18.5
$ws.Range("AB5").Value = 17
$ws.Range("AC5").Value = 8.199999999999999
$ws.Range("AD5").Value = 9.800000000000001
$ws.Range("AE5").Value = 19.5
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 20
$ws.Range("AH5").Value = 22
$ws.Range("AI5").Value = 38
$ws.Range("AJ5").Value = 130
$ws.Range("AK5").Value = 75
$ws.Range("AL5").Value = 75
$ws.Range("AM5").Value = 130
$ws.Range("AN5").Value = 90
$ws.Range("AO5").Value = 12.5
$ws.Range("F6").Value = 2.78
$ws.Range("G6").Value = 2.84
$ws.Range("H6").Value = 2.72
$ws.Range("I6").Value = 2.78
$ws.Range("J6").Value = 3.55
$ws.Range("K6").Value = 3.6
$ws.Range("L6").Value = 1.4
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 3.65
$ws.Range("O6").Value = 1.32
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.34
$ws.Range("S6").Value = 3.35
$ws.Range("T6").Value = 1.73
$ws.Range("U6").Value = 2.16
$ws.Range("V6").Value = 1.56
$ws.Range("W6").Value = 1.54
$ws.Range("X6").Value = 15
$ws.Range("Y6").Value = 12
$ws.Range("Z6").Value = 19
$ws.Range("AA6").Value = 44
$ws.Range("AB6").Value = 12
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 32
$ws.Range("AF6").Value = 18.5
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 17.5
$ws.Range("AI6").Value = 44
$ws.Range("AJ6").Value = 44
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 44
$ws.Range("AM6").Value = 95
$ws.Range("AN6").Value = 29
$ws.Range("AO6").Value = 28
$ws.Range("L7").Value = 1.35
$ws.Range("P7").Value = 1.95
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.36
$ws.Range("T7").Value = 1.79
$ws.Range("V7").Value = 1.7
$ws.Range("W7").Value = 1.43
$ws.Range("X7").Value = 14.5
$ws.Range("Z7").Value = 14.5
$ws.Range("AA7").Value = 32
$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 25
$ws.Range("AF7").Value = 22
$ws.Range("AG7").Value = 14
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 40
$ws.Range("AJ7").Value = 55
$ws.Range("AK7").Value = 38
$ws.Range("AN7").Value = 34
$ws.Range("AO7").Value = 20
$ws.Range("F8").Value = 1.65
$ws.Range("G8").Value = 1.67
$ws.Range("H8").Value = 5.8
$ws.Range("K8").Value = 4.5
$ws.Range("L8").Value = 1.33
$ws.Range("N8").Value = 4.3
$ws.Range("Q8").Value = 1.79
$ws.Range("S8").Value = 2.88
$ws.Range("T8").Value = 1.79
$ws.Range("U8").Value = 2.04
$ws.Range("V8").Value = 1.19
$ws.Range("W8").Value = 2.48
$ws.Range("X8").Value = 17.5
$ws.Range("Z8").Value = 50
$ws.Range("AA8").Value = 170
$ws.Range("AB8").Value = 10
$ws.Range("AC8").Value = 9.800000000000001
$ws.Range("AD8").Value = 23
$ws.Range("AF8").Value = 10.5
$ws.Range("AG8").Value = 10.5
$ws.Range("AJ8").Value = 16
$ws.Range("AK8").Value = 20
$ws.Range("AL8").Value = 34
$ws.Range("AM8").Value = 110
$ws.Range("AO8").Value = 90
